$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 4).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 99 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $val = $cell.Value2
    if ($val -eq "MOLLY MCNINCH") {
        $cell.Value = "T"
    } elseif ($val -eq "STUDENT") {
        $cell.Value = "S"
    }
}
